$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 142
$ws.Range("F4").Value = 801
$ws.Range("F6").Value = 644
$ws.Range("F7").Value = 1213
$ws.Range("F9").Value = 792
$ws.Range("F10").Value = 682
$ws.Range("F13").Value = 353
$ws.Range("F15").Value = 898
$ws.Range("F16").Value = 9727
$ws.Range("F17").Value = 600
$ws.Range("F19").Value = 312
$ws.Range("F22").Value = 259
$ws.Range("F24").Value = 24
$ws.Range("F25").Value = 290
$ws.Range("F26").Value = 476
$ws.Range("F27").Value = 177
$ws.Range("F29").Value = 258
$ws.Range("F31").Value = 262
$ws.Range("F32").Value = 59
$ws.Range("F34").Value = 19
$ws.Range("F36").Value = 185
$ws.Range("F38").Value = 36
$ws.Range("F39").Value = 91

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 15
$ws.Range("F7").Value = 121
$ws.Range("F11").Value = 4439
$ws.Range("F12").Value = 82
$ws.Range("F16").Value = 264
$ws.Range("F20").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 808

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 808
$ws.Range("F5").Value = 142
$ws.Range("F6").Value = 801
$ws.Range("F7").Value = 15
$ws.Range("F10").Value = 644
$ws.Range("F11").Value = 1213
$ws.Range("F13").Value = 121
$ws.Range("F14").Value = 792
$ws.Range("F15").Value = 682
$ws.Range("F17").Value = 353
$ws.Range("F19").Value = 898
$ws.Range("F20").Value = 9727
$ws.Range("F22").Value = 600
$ws.Range("F25").Value = 259
$ws.Range("F27").Value = 24
$ws.Range("F28").Value = 476
$ws.Range("F29").Value = 177
$ws.Range("F30").Value = 82
$ws.Range("F31").Value = 82
$ws.Range("F36").Value = 258
$ws.Range("F38").Value = 262
$ws.Range("F39").Value = 59
$ws.Range("F41").Value = 19
$ws.Range("F45").Value = 11
$ws.Range("F46").Value = 185
